$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(121, 8).Value = 1454.5
$ws.Cells.Item(121, 9).Value = 575
$ws.Cells.Item(121, 10).Value = 1630.4
$ws.Cells.Item(121, 11).Value = 1725
$ws.Cells.Item(121, 12).Value = 4891.200000000001
$ws.Cells.Item(121, 13).Value = 22
$ws.Cells.Item(121, 14).Value = -8385.200000000001
$ws.Cells.Item(127, 8).Value = 3848
$ws.Cells.Item(127, 9).Value = 2676
$ws.Cells.Item(127, 10).Value = 7950
$ws.Cells.Item(127, 11).Value = 8028
$ws.Cells.Item(127, 12).Value = 23850
$ws.Cells.Item(127, 13).Value = -3068
$ws.Cells.Item(127, 14).Value = -33770
$ws.Cells.Item(129, 8).Value = 1038.2759
$ws.Cells.Item(129, 9).Value = 1042
$ws.Cells.Item(129, 10).Value = 1037.5
$ws.Cells.Item(129, 11).Value = 3126
$ws.Cells.Item(129, 12).Value = 3112.5
$ws.Cells.Item(129, 13).Value = 1874
$ws.Cells.Item(129, 14).Value = -13112.5
$ws.Cells.Item(132, 8).Value = 979.75555
$ws.Cells.Item(132, 9).Value = 977.0227
$ws.Cells.Item(132, 10).Value = 1100
$ws.Cells.Item(132, 11).Value = 2931.0681
$ws.Cells.Item(132, 12).Value = 3300
$ws.Cells.Item(132, 13).Value = -401.0681
$ws.Cells.Item(132, 14).Value = -8360
$ws.Cells.Item(137, 8).Value = 2304
$ws.Cells.Item(137, 9).Value = 1998
$ws.Cells.Item(137, 10).Value = 2387.4546
$ws.Cells.Item(137, 11).Value = 5994
$ws.Cells.Item(137, 12).Value = 7162.3638
$ws.Cells.Item(137, 13).Value = -3444
$ws.Cells.Item(137, 14).Value = -12262.3638
$ws.Cells.Item(138, 8).Value = 6028.5293
$ws.Cells.Item(138, 9).Value = 5415.3335
$ws.Cells.Item(138, 10).Value = 6217.205
$ws.Cells.Item(138, 11).Value = 16246.0005
$ws.Cells.Item(138, 12).Value = 18651.615
$ws.Cells.Item(138, 13).Value = -11106.0005
$ws.Cells.Item(138, 14).Value = -28931.615

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3075.3452
$ws.Cells.Item(32, 9).Value = 2505.9578
$ws.Cells.Item(32, 10).Value = 6185.077
$ws.Cells.Item(32, 11).Value = 2505.9578
$ws.Cells.Item(32, 12).Value = 6185.077
$ws.Cells.Item(32, 13).Value = -2218.9578
$ws.Cells.Item(45, 8).Value = 1334.2
$ws.Cells.Item(45, 9).Value = 1185.4286
$ws.Cells.Item(45, 10).Value = 1681.3334
$ws.Cells.Item(45, 11).Value = 1185.4286
$ws.Cells.Item(45, 12).Value = 1681.3334
$ws.Cells.Item(45, 13).Value = -808.4286
$ws.Cells.Item(45, 14).Value = -2435.3334
$ws.Cells.Item(88, 8).Value = 2955.5
$ws.Cells.Item(88, 9).Value = 2077
$ws.Cells.Item(88, 10).Value = 3638.7778
$ws.Cells.Item(88, 11).Value = 2077
$ws.Cells.Item(88, 12).Value = 3638.7778
$ws.Cells.Item(88, 13).Value = -1671
$ws.Cells.Item(88, 14).Value = -4450.7778
$ws.Cells.Item(91, 8).Value = 2955.5
$ws.Cells.Item(91, 9).Value = 2077
$ws.Cells.Item(91, 10).Value = 3638.7778
$ws.Cells.Item(91, 11).Value = 2077
$ws.Cells.Item(91, 12).Value = 3638.7778
$ws.Cells.Item(91, 13).Value = -673
$ws.Cells.Item(91, 14).Value = -6446.7778

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 7423.7334
$ws.Cells.Item(134, 9).Value = 8244.25
$ws.Cells.Item(134, 10).Value = 4141.6665
$ws.Cells.Item(134, 11).Value = 24732.75
$ws.Cells.Item(134, 12).Value = 12424.9995
$ws.Cells.Item(134, 13).Value = -22197.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 967.875
$ws.Cells.Item(16, 9).Value = 957.1667
$ws.Cells.Item(16, 10).Value = 1000
$ws.Cells.Item(16, 11).Value = 957.1667
$ws.Cells.Item(16, 12).Value = 1000
$ws.Cells.Item(16, 13).Value = -670.1667
$ws.Cells.Item(16, 14).Value = -1574
$ws.Cells.Item(31, 8).Value = 1251.1482
$ws.Cells.Item(31, 9).Value = 904.4545000000001
$ws.Cells.Item(31, 10).Value = 1489.5
$ws.Cells.Item(31, 11).Value = 904.4545000000001
$ws.Cells.Item(31, 12).Value = 1489.5
$ws.Cells.Item(31, 13).Value = -609.4545000000001
$ws.Cells.Item(31, 14).Value = -2079.5
$ws.Cells.Item(34, 8).Value = 1251.1482
$ws.Cells.Item(34, 9).Value = 904.4545000000001
$ws.Cells.Item(34, 10).Value = 1489.5
$ws.Cells.Item(34, 11).Value = 904.4545000000001
$ws.Cells.Item(34, 12).Value = 1489.5
$ws.Cells.Item(34, 13).Value = -702.4545000000001
$ws.Cells.Item(34, 14).Value = -1893.5
$ws.Cells.Item(113, 8).Value = 967.875
$ws.Cells.Item(113, 9).Value = 957.1667
$ws.Cells.Item(113, 10).Value = 1000
$ws.Cells.Item(113, 11).Value = 957.1667
$ws.Cells.Item(113, 12).Value = 1000
$ws.Cells.Item(113, 13).Value = 1212.8333
$ws.Cells.Item(113, 14).Value = -5340
$ws.Cells.Item(129, 7).Value = 35378
$ws.Cells.Item(129, 8).Value = 30000
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 30000
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 30000
$ws.Cells.Item(129, 14).Value = -40000
$ws.Cells.Item(130, 7).Value = 34689
$ws.Cells.Item(130, 8).Value = 23800
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 23800
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 23800
$ws.Cells.Item(130, 14).Value = -33840
$ws.Cells.Item(131, 7).Value = 35461
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(132, 7).Value = 44019
$ws.Cells.Item(132, 8).Value = 2453.5
$ws.Cells.Item(132, 9).Value = 1172.25
$ws.Cells.Item(132, 10).Value = 3991
$ws.Cells.Item(132, 11).Value = 3516.75
$ws.Cells.Item(132, 12).Value = 11973
$ws.Cells.Item(132, 13).Value = -986.75
$ws.Cells.Item(132, 14).Value = -17033
$ws.Cells.Item(133, 7).Value = 43328
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(134, 7).Value = 44020
$ws.Cells.Item(134, 8).Value = 2135.5
$ws.Cells.Item(134, 9).Value = 2135.5
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 6406.5
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -3871.5
$ws.Cells.Item(135, 7).Value = 42008
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(137, 7).Value = 43231
$ws.Cells.Item(137, 8).Value = 60780
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 60780
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 60780
$ws.Cells.Item(137, 14).Value = -70980
$ws.Cells.Item(138, 7).Value = 42302
$ws.Cells.Item(138, 8).Value = 24500
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 24500
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 24500
$ws.Cells.Item(138, 14).Value = -34780
$ws.Cells.Item(139, 7).Value = 43258
$ws.Cells.Item(139, 8).Value = 51923
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 51923
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 51923
$ws.Cells.Item(139, 14).Value = -62203
$ws.Cells.Item(140, 7).Value = 42455
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(141, 7).Value = 43345
$ws.Cells.Item(141, 8).Value = 30882
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 30882
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 30882
$ws.Cells.Item(141, 14).Value = -41242

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 3366
$ws.Cells.Item(68, 9).Value = 1399
$ws.Cells.Item(68, 10).Value = 3967.0278
$ws.Cells.Item(68, 11).Value = 4197
$ws.Cells.Item(68, 12).Value = 11901.0834
$ws.Cells.Item(68, 13).Value = -3386
$ws.Cells.Item(68, 14).Value = -13523.0834
$ws.Cells.Item(71, 8).Value = 3366
$ws.Cells.Item(71, 9).Value = 1399
$ws.Cells.Item(71, 10).Value = 3967.0278
$ws.Cells.Item(71, 11).Value = 12591
$ws.Cells.Item(71, 12).Value = 35703.25019999999
$ws.Cells.Item(71, 13).Value = -8535
$ws.Cells.Item(71, 14).Value = -43815.25019999999
$ws.Cells.Item(131, 8).Value = 12213534
$ws.Cells.Item(131, 9).Value = 41667136
$ws.Cells.Item(131, 10).Value = 25836.414
$ws.Cells.Item(131, 11).Value = 125001408
$ws.Cells.Item(131, 12).Value = 77509.242
$ws.Cells.Item(131, 13).Value = -124996368
$ws.Cells.Item(131, 14).Value = -87589.242

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 17773.834
$ws.Cells.Item(16, 9).Value = 20528.6
$ws.Cells.Item(16, 10).Value = 4000
$ws.Cells.Item(16, 11).Value = 20528.6
$ws.Cells.Item(16, 12).Value = 4000
$ws.Cells.Item(16, 13).Value = -20358.6
$ws.Cells.Item(16, 14).Value = -4340
$ws.Cells.Item(55, 8).Value = 306.84616
$ws.Cells.Item(55, 9).Value = 248.9
$ws.Cells.Item(55, 10).Value = 500
$ws.Cells.Item(55, 11).Value = 248.9
$ws.Cells.Item(55, 12).Value = 500
$ws.Cells.Item(55, 13).Value = -75.90000000000001
$ws.Cells.Item(55, 14).Value = -846
$ws.Cells.Item(82, 8).Value = 4696.6665
$ws.Cells.Item(82, 9).Value = 2000
$ws.Cells.Item(82, 10).Value = 5236
$ws.Cells.Item(82, 11).Value = 2000
$ws.Cells.Item(82, 12).Value = 5236
$ws.Cells.Item(82, 13).Value = -1639
$ws.Cells.Item(82, 14).Value = -5958
$ws.Cells.Item(85, 8).Value = 4696.6665
$ws.Cells.Item(85, 9).Value = 2000
$ws.Cells.Item(85, 10).Value = 5236
$ws.Cells.Item(85, 11).Value = 2000
$ws.Cells.Item(85, 12).Value = 5236
$ws.Cells.Item(85, 13).Value = -752
$ws.Cells.Item(85, 14).Value = -7732
$ws.Cells.Item(132, 8).Value = 2693.0176
$ws.Cells.Item(132, 9).Value = 1306
$ws.Cells.Item(132, 10).Value = 3856.3225
$ws.Cells.Item(132, 11).Value = 3918
$ws.Cells.Item(132, 12).Value = 11568.9675
$ws.Cells.Item(132, 13).Value = -1388
$ws.Cells.Item(136, 8).Value = 3759.9565
$ws.Cells.Item(136, 9).Value = 2809.394
$ws.Cells.Item(136, 10).Value = 6172.923
$ws.Cells.Item(136, 11).Value = 8428.181999999999
$ws.Cells.Item(136, 12).Value = 18518.769
$ws.Cells.Item(136, 13).Value = -5878.181999999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3000
$ws.Cells.Item(81, 9).Value = 3000
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 6000
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = -4939
$ws.Cells.Item(84, 8).Value = 3000
$ws.Cells.Item(84, 9).Value = 3000
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 30000
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = -24696
$ws.Cells.Item(126, 8).Value = 9987.951999999999
$ws.Cells.Item(126, 9).Value = 9985.944
$ws.Cells.Item(126, 10).Value = 10000
$ws.Cells.Item(126, 11).Value = 29957.832
$ws.Cells.Item(126, 12).Value = 30000
$ws.Cells.Item(126, 13).Value = -27487.832
$ws.Cells.Item(126, 14).Value = -34940
$ws.Cells.Item(132, 8).Value = 1221.9166
$ws.Cells.Item(132, 9).Value = 976.2727
$ws.Cells.Item(132, 10).Value = 1429.7693
$ws.Cells.Item(132, 11).Value = 2928.8181
$ws.Cells.Item(132, 12).Value = 4289.3079
$ws.Cells.Item(132, 13).Value = -398.8181
$ws.Cells.Item(132, 14).Value = -9349.3079
$ws.Cells.Item(81, 14).ClearContents()
$ws.Cells.Item(84, 14).ClearContents()
